# Update countries & provincias Spain
# - Reorders Haiti, Santa Lucia and San Bartolome within the country list
#   (their shared-string position moves, shifting which row shows which
#   country and which numbers go with it)
# - Refreshes the "Datos actualizados" timestamp
# - Refreshes case counts for China, Kazajistan and Honduras

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 05:05"

# --- China (row 19): new case counts ---
$ws.Range("B19").Value = 82999
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 78302
$ws.Range("E19").Value = 63

# --- Kazajistan (row 54): new case counts ---
$ws.Range("B54").Value = 10382
$ws.Range("C54").Value = 450
$ws.Range("E54").Value = 5288

# --- Honduras (row 71): new case counts ---
$ws.Range("B71").Value = 4886
$ws.Range("C71").Value = 134
$ws.Range("D71").Value = 528
$ws.Range("E71").Value = 4159
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 199

# --- Haiti is inserted into the list before Sri Lanka, pushing Sri
#     Lanka, Eslovaquia, Nueva Zelanda and Eslovenia down one row each
#     (rows 100-104) ---
$ws.Range("A100").Value = "Haiti"
$ws.Range("B100").Value = 1584
$ws.Range("C100").Value = 141
$ws.Range("D100").Value = 22
$ws.Range("E100").Value = 1527
$ws.Range("H100").Value = 35

$ws.Range("A101").Value = "Sri Lanka"
$ws.Range("B101").Value = 1558
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 754
$ws.Range("E101").Value = 794
$ws.Range("H101").Value = 10

$ws.Range("A102").Value = "Eslovaquia"
$ws.Range("B102").Value = 1520
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 1338
$ws.Range("E102").Value = 154
$ws.Range("H102").Value = 28

$ws.Range("A103").Value = "Nueva Zelanda"
$ws.Range("B103").Value = 1504
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 1481
$ws.Range("E103").Value = 1
$ws.Range("H103").Value = 22

$ws.Range("A104").Value = "Eslovenia"
$ws.Range("B104").Value = 1473
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 1357
$ws.Range("E104").Value = 8
$ws.Range("H104").Value = 108

# --- Santa Lucia moves ahead of Belice (rows 200-201 swap identity &
#     their Casos activos / Muertes figures) ---
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# --- San Bartolome moves ahead of "Bonaire, San Eustaquio y Saba"
#     (rows 215-216 swap identity; underlying figures are identical so
#     no numeric change is required) ---
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
